$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 574; everything from 574 downward shifts down by one.
$ws.Rows.Item(574).Insert()

# Populate the newly inserted row 574 with the new daily entry.
# Use a Text number format while assigning the date-like string so it is
# stored as literal text (matching the other date cells) instead of being
# auto-parsed into a date serial number, then clear the format back to the
# sheet's default so no stray style is left behind.
$ws.Range("A574").NumberFormat = "@"
$ws.Range("A574").Value = "2026/01/07"
$ws.Range("A574").ClearFormats()

$ws.Range("B574").Value = "水"
$ws.Range("C574").Value = 16
$ws.Range("D574").Value = 24
